$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 298.77777
$ws.Range("I2").Value = 277.64285
$ws.Range("J2").Value = 372.75
$ws.Range("K2").Value = 277.64285
$ws.Range("L2").Value = 372.75
$ws.Range("M2").Value = -164.64285
$ws.Range("N2").Value = -598.75
$ws.Range("H4").Value = 296.42856
$ws.Range("I4").Value = 236.36363
$ws.Range("J4").Value = 516.6667
$ws.Range("K4").Value = 236.36363
$ws.Range("L4").Value = 516.6667
$ws.Range("M4").Value = -122.36363
$ws.Range("N4").Value = -744.6667
$ws.Range("H8").Value = 1300.381
$ws.Range("I8").Value = 1021.2857
$ws.Range("J8").Value = 1858.5714
$ws.Range("K8").Value = 3063.8571
$ws.Range("L8").Value = 5575.7142
$ws.Range("M8").Value = -2924.8571
$ws.Range("N8").Value = -5853.7142
$ws.Range("H9").Value = 5026.2856
$ws.Range("I9").Value = 193.33333
$ws.Range("K9").Value = 193.33333
$ws.Range("M9").Value = -24.33332999999999
$ws.Range("H15").Value = 1905.8
$ws.Range("I15").Value = 1905.8
$ws.Range("K15").Value = 5717.4
$ws.Range("M15").Value = -5548.4
$ws.Range("H21").Value = 32500
$ws.Range("J21").Value = 15000
$ws.Range("L21").Value = 15000
$ws.Range("N21").Value = -15936
$ws.Range("H23").Value = 32500
$ws.Range("J23").Value = 15000
$ws.Range("L23").Value = 15000
$ws.Range("N23").Value = -15468
$ws.Range("H29").Value = 3591.158
$ws.Range("I29").Value = 2178
$ws.Range("J29").Value = 4618.909
$ws.Range("K29").Value = 6534
$ws.Range("L29").Value = 13856.727
$ws.Range("M29").Value = -6253
$ws.Range("N29").Value = -14418.727
$ws.Range("H31").Value = 1030.7693
$ws.Range("I31").Value = 763.63635
$ws.Range("J31").Value = 2500
$ws.Range("K31").Value = 2290.90905
$ws.Range("L31").Value = 7500
$ws.Range("M31").Value = -2060.90905
$ws.Range("N31").Value = -7960
$ws.Range("H32").Value = 6013.3335
$ws.Range("J32").Value = 7656.4
$ws.Range("L32").Value = 7656.4
$ws.Range("N32").Value = -8308.4
$ws.Range("H33").Value = 400.8125
$ws.Range("I33").Value = 336.3846
$ws.Range("K33").Value = 336.3846
$ws.Range("M33").Value = -107.3846
$ws.Range("H38").Value = 190.75
$ws.Range("I38").Value = 26.272728
$ws.Range("J38").Value = 2000
$ws.Range("K38").Value = 78.818184
$ws.Range("L38").Value = 6000
$ws.Range("M38").Value = 293.181816
$ws.Range("N38").Value = -6744
$ws.Range("H40").Value = 4150.033
$ws.Range("J40").Value = 6162.75
$ws.Range("L40").Value = 6162.75
$ws.Range("N40").Value = -6512.75
$ws.Range("H41").Value = 263.125
$ws.Range("I41").Value = 98.333336
$ws.Range("J41").Value = 362
$ws.Range("K41").Value = 98.333336
$ws.Range("L41").Value = 362
$ws.Range("M41").Value = 341.666664
$ws.Range("N41").Value = -1242
$ws.Range("H55").Value = 477.1613
$ws.Range("I55").Value = 356.70587
$ws.Range("J55").Value = 623.4286
$ws.Range("K55").Value = 356.70587
$ws.Range("L55").Value = 623.4286
$ws.Range("M55").Value = -142.70587
$ws.Range("N55").Value = -1051.4286
$ws.Range("H58").Value = 1846.7142
$ws.Range("I58").Value = 1472.8889
$ws.Range("J58").Value = 2519.6
$ws.Range("K58").Value = 4418.6667
$ws.Range("L58").Value = 7558.799999999999
$ws.Range("M58").Value = -4268.6667
$ws.Range("N58").Value = -7858.799999999999
$ws.Range("H69").Value = 7431.25
$ws.Range("J69").Value = 7650.722
$ws.Range("L69").Value = 22952.166
$ws.Range("N69").Value = -24700.166
$ws.Range("H70").Value = 56474.55
$ws.Range("I70").Value = 116054.555
$ws.Range("J70").Value = 7727.273
$ws.Range("K70").Value = 348163.665
$ws.Range("L70").Value = 23181.819
$ws.Range("M70").Value = -347893.665
$ws.Range("N70").Value = -23721.819
$ws.Range("H72").Value = 7431.25
$ws.Range("J72").Value = 7650.722
$ws.Range("L72").Value = 68856.49799999999
$ws.Range("N72").Value = -77592.49799999999
$ws.Range("H73").Value = 56474.55
$ws.Range("I73").Value = 116054.555
$ws.Range("J73").Value = 7727.273
$ws.Range("K73").Value = 348163.665
$ws.Range("L73").Value = 23181.819
$ws.Range("M73").Value = -347227.665
$ws.Range("N73").Value = -25053.819
$ws.Range("H86").Value = 4744.4814
$ws.Range("I86").Value = 3828.8572
$ws.Range("K86").Value = 3828.8572
$ws.Range("M86").Value = -2705.8572
$ws.Range("H89").Value = 4744.4814
$ws.Range("I89").Value = 3828.8572
$ws.Range("K89").Value = 19144.286
$ws.Range("M89").Value = -13528.286
$ws.Range("H98").Value = 306147.75
$ws.Range("I98").Value = 1495.8572
$ws.Range("J98").Value = 1017002.2
$ws.Range("K98").Value = 1495.8572
$ws.Range("L98").Value = 1017002.2
$ws.Range("M98").Value = 2.142800000000079
$ws.Range("N98").Value = -1019998.2
$ws.Range("H112").Value = 1313.3549
$ws.Range("J112").Value = 1338.3704
$ws.Range("L112").Value = 4015.1112
$ws.Range("N112").Value = -6231.1112
$ws.Range("H113").Value = 5473.357
$ws.Range("I113").Value = 4074.8572
$ws.Range("J113").Value = 6871.857
$ws.Range("K113").Value = 4074.8572
$ws.Range("L113").Value = 6871.857
$ws.Range("M113").Value = -820.8571999999999
$ws.Range("N113").Value = -13379.857
$ws.Range("H116").Value = 3920
$ws.Range("I116").Value = 3305.75
$ws.Range("K116").Value = 3305.75
$ws.Range("M116").Value = 136.25
$ws.Range("H122").Value = 306147.75
$ws.Range("I122").Value = 1495.8572
$ws.Range("J122").Value = 1017002.2
$ws.Range("K122").Value = 4487.571599999999
$ws.Range("L122").Value = 3051006.6
$ws.Range("M122").Value = -2037.571599999999
$ws.Range("N122").Value = -3055906.6
$ws.Range("H127").Value = 389
$ws.Range("I127").Value = 389
$ws.Range("K127").Value = 1167
$ws.Range("M127").Value = 3793
$ws.Range("H135").Value = 1406.4706
$ws.Range("I135").Value = 925
$ws.Range("J135").Value = 5017.5
$ws.Range("K135").Value = 8325
$ws.Range("L135").Value = 45157.5
$ws.Range("M135").Value = -5790
$ws.Range("N135").Value = -50227.5
$ws.Range("H137").Value = 3011.9138
$ws.Range("I137").Value = 2013.5
$ws.Range("K137").Value = 6040.5
$ws.Range("M137").Value = -3490.5
$ws.Range("H138").Value = 2647.4197
$ws.Range("I138").Value = 1210.6451
$ws.Range("J138").Value = 3538.22
$ws.Range("K138").Value = 3631.9353
$ws.Range("L138").Value = 10614.66
$ws.Range("M138").Value = 1508.0647
$ws.Range("N138").Value = -20894.66
$ws.Range("H141").Value = 3474
$ws.Range("I141").Value = 3474
$ws.Range("K141").Value = 10422
$ws.Range("M141").Value = -5242

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 171.33333
$ws.Range("H32").Value = 2536.5334
$ws.Range("I32").Value = 2059.0417
$ws.Range("J32").Value = 13996.333
$ws.Range("K32").Value = 2059.0417
$ws.Range("L32").Value = 13996.333
$ws.Range("M32").Value = -1772.0417
$ws.Range("N32").Value = -14570.333
$ws.Range("H45").Value = 4739
$ws.Range("I45").Value = 4800
$ws.Range("K45").Value = 4800
$ws.Range("M45").Value = -4423
$ws.Range("H60").Value = 64833.168
$ws.Range("I60").Value = 64833.168
$ws.Range("K60").Value = 64833.168
$ws.Range("M60").Value = -64100.168
$ws.Range("H61").Value = 4285.0605
$ws.Range("I61").Value = 3078.3928
$ws.Range("K61").Value = 3078.3928
$ws.Range("M61").Value = -2866.3928
$ws.Range("H74").Value = 30307570
$ws.Range("I74").Value = 33336424
$ws.Range("K74").Value = 33336424
$ws.Range("M74").Value = -33335550
$ws.Range("H77").Value = 30307570
$ws.Range("I77").Value = 33336424
$ws.Range("K77").Value = 166682120
$ws.Range("M77").Value = -166677752
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H97").Value = 1157.9333
$ws.Range("I97").Value = 807.1111
$ws.Range("J97").Value = 1684.1666
$ws.Range("K97").Value = 807.1111
$ws.Range("L97").Value = 1684.1666
$ws.Range("M97").Value = -311.1111
$ws.Range("N97").Value = -2676.1666
$ws.Range("H109").Value = 62039.2
$ws.Range("J109").Value = 62039.2
$ws.Range("L109").Value = 62039.2
$ws.Range("N109").Value = -64813.2
$ws.Range("H122").Value = 4265.2354
$ws.Range("I122").Value = 4282.5
$ws.Range("J122").Value = 4249.8887
$ws.Range("K122").Value = 12847.5
$ws.Range("L122").Value = 12749.6661
$ws.Range("M122").Value = -10397.5
$ws.Range("N122").Value = -17649.6661
$ws.Range("H132").Value = 5134.125
$ws.Range("I132").Value = 3009.8333
$ws.Range("J132").Value = 11507
$ws.Range("K132").Value = 9029.499899999999
$ws.Range("L132").Value = 34521
$ws.Range("M132").Value = -6499.499899999999
$ws.Range("N132").Value = -39581
$ws.Range("H136").Value = 4285.0605
$ws.Range("I136").Value = 3078.3928
$ws.Range("K136").Value = 9235.178400000001
$ws.Range("M136").Value = -6685.178400000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2272.158
$ws.Range("I94").Value = 2257.5334
$ws.Range("J94").Value = 2327
$ws.Range("K94").Value = 2257.5334
$ws.Range("L94").Value = 2327
$ws.Range("M94").Value = -1806.5334
$ws.Range("N94").Value = -3229
$ws.Range("H96").Value = 12312.8
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()
$ws.Range("H105").Value = 14647.682
$ws.Range("I105").Value = 16863.924
$ws.Range("K105").Value = 16863.924
$ws.Range("M105").Value = -15116.924
$ws.Range("H134").Value = 2772.4211
$ws.Range("I134").Value = 2772.4211
$ws.Range("K134").Value = 8317.263300000001
$ws.Range("M134").Value = -5782.263300000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 1006200
$ws.Range("J4").Value = 1257500
$ws.Range("L4").Value = 1257500
$ws.Range("N4").Value = -1257724
$ws.Range("H16").Value = 1458.1818
$ws.Range("I16").Value = 553.375
$ws.Range("J16").Value = 3871
$ws.Range("K16").Value = 553.375
$ws.Range("L16").Value = 3871
$ws.Range("M16").Value = -266.375
$ws.Range("N16").Value = -4445
$ws.Range("H31").Value = 176654.67
$ws.Range("J31").Value = 176654.67
$ws.Range("L31").Value = 176654.67
$ws.Range("N31").Value = -177244.67
$ws.Range("H34").Value = 176654.67
$ws.Range("J34").Value = 176654.67
$ws.Range("L34").Value = 176654.67
$ws.Range("N34").Value = -177058.67
$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("M55").ClearContents()
$ws.Range("N55").ClearContents()
$ws.Range("H58").Value = 12812.375
$ws.Range("I58").Value = 3501
$ws.Range("J58").Value = 14142.571
$ws.Range("K58").Value = 3501
$ws.Range("L58").Value = 14142.571
$ws.Range("M58").Value = -3298
$ws.Range("N58").Value = -14548.571
$ws.Range("H59").Value = 74997.664
$ws.Range("J59").Value = 74997.664
$ws.Range("L59").Value = 74997.664
$ws.Range("N59").Value = -77287.664
$ws.Range("H62").Value = 10550.5
$ws.Range("I62").Value = 4099.5
$ws.Range("J62").Value = 17001.5
$ws.Range("K62").Value = 4099.5
$ws.Range("L62").Value = 17001.5
$ws.Range("M62").Value = -3475.5
$ws.Range("N62").Value = -18249.5
$ws.Range("H65").Value = 10550.5
$ws.Range("I65").Value = 4099.5
$ws.Range("J65").Value = 17001.5
$ws.Range("K65").Value = 20497.5
$ws.Range("L65").Value = 85007.5
$ws.Range("M65").Value = -17377.5
$ws.Range("N65").Value = -91247.5
$ws.Range("H86").Value = 8625
$ws.Range("I86").Value = 8166.6665
$ws.Range("J86").Value = 10000
$ws.Range("K86").Value = 8166.6665
$ws.Range("L86").Value = 10000
$ws.Range("M86").Value = -7043.6665
$ws.Range("N86").Value = -12246
$ws.Range("H89").Value = 8625
$ws.Range("I89").Value = 8166.6665
$ws.Range("J89").Value = 10000
$ws.Range("K89").Value = 40833.3325
$ws.Range("L89").Value = 50000
$ws.Range("M89").Value = -35217.3325
$ws.Range("N89").Value = -61232
$ws.Range("H95").Value = 17812
$ws.Range("J95").Value = 17812
$ws.Range("L95").Value = 17812
$ws.Range("N95").Value = -23304
$ws.Range("H113").Value = 1458.1818
$ws.Range("I113").Value = 553.375
$ws.Range("J113").Value = 3871
$ws.Range("K113").Value = 553.375
$ws.Range("L113").Value = 3871
$ws.Range("M113").Value = 1616.625
$ws.Range("N113").Value = -8211
$ws.Range("H132").Value = 3572.28
$ws.Range("I132").Value = 2231
$ws.Range("K132").Value = 6693
$ws.Range("M132").Value = -4163
$ws.Range("H134").Value = 2902.4443
$ws.Range("I134").Value = 2016.7858
$ws.Range("J134").Value = 6002.25
$ws.Range("K134").Value = 6050.357400000001
$ws.Range("L134").Value = 18006.75
$ws.Range("M134").Value = -3515.357400000001
$ws.Range("N134").Value = -23076.75
$ws.Range("H136").Value = 12812.375
$ws.Range("I136").Value = 3501
$ws.Range("J136").Value = 14142.571
$ws.Range("K136").Value = 10503
$ws.Range("L136").Value = 42427.713
$ws.Range("M136").Value = -7953
$ws.Range("N136").Value = -47527.713

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 157.33333
$ws.Range("J2").Value = 181.8
$ws.Range("L2").Value = 1090.8
$ws.Range("N2").Value = -1316.8
$ws.Range("H5").Value = 2775
$ws.Range("J5").Value = 4466.8335
$ws.Range("L5").Value = 13400.5005
$ws.Range("N5").Value = -13624.5005
$ws.Range("H7").Value = 411.75
$ws.Range("I7").Value = 444.83334
$ws.Range("J7").Value = 312.5
$ws.Range("K7").Value = 1334.50002
$ws.Range("L7").Value = 937.5
$ws.Range("M7").Value = -1222.50002
$ws.Range("N7").Value = -1161.5
$ws.Range("H12").Value = 8.800000000000001
$ws.Range("J12").Value = 8.800000000000001
$ws.Range("L12").Value = 26.4
$ws.Range("N12").Value = -372.4
$ws.Range("H17").Value = 575.75
$ws.Range("I17").Value = 228.77777
$ws.Range("J17").Value = 1616.6666
$ws.Range("K17").Value = 686.33331
$ws.Range("L17").Value = 4849.9998
$ws.Range("M17").Value = -517.33331
$ws.Range("N17").Value = -5187.9998
$ws.Range("H37").Value = 212375.44
$ws.Range("J37").Value = 212375.44
$ws.Range("L37").Value = 637126.3200000001
$ws.Range("N37").Value = -637350.3200000001
$ws.Range("H38").Value = 46
$ws.Range("I38").Value = 52.833332
$ws.Range("J38").Value = 37.8
$ws.Range("K38").Value = 158.499996
$ws.Range("L38").Value = 113.4
$ws.Range("M38").Value = 188.500004
$ws.Range("N38").Value = -807.4
$ws.Range("H46").Value = 2632.25
$ws.Range("I46").Value = 186.83333
$ws.Range("J46").Value = 4099.5
$ws.Range("K46").Value = 560.49999
$ws.Range("L46").Value = 12298.5
$ws.Range("M46").Value = -469.49999
$ws.Range("N46").Value = -12480.5
$ws.Range("H50").Value = 25647242
$ws.Range("I50").Value = 47620616
$ws.Range("K50").Value = 142861848
$ws.Range("M50").Value = -142861367
$ws.Range("H53").Value = 25647242
$ws.Range("I53").Value = 47620616
$ws.Range("K53").Value = 142861848
$ws.Range("M53").Value = -142861367
$ws.Range("H135").Value = 2775
$ws.Range("J135").Value = 4466.8335
$ws.Range("L135").Value = 40201.5015
$ws.Range("N135").Value = -45271.5015
$ws.Range("H137").Value = 91813.09
$ws.Range("J137").Value = 501500.5
$ws.Range("L137").Value = 1504501.5
$ws.Range("N137").Value = -1514701.5
$ws.Range("H139").Value = 4103.8667
$ws.Range("I139").Value = 2321.7778
$ws.Range("K139").Value = 6965.3334
$ws.Range("M139").Value = -1825.3334

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 46665.5
$ws.Range("J15").Value = 46665.5
$ws.Range("L15").Value = 46665.5
$ws.Range("N15").Value = -47241.5
$ws.Range("H70").Value = 12038.723
$ws.Range("I70").Value = 11957.429
$ws.Range("K70").Value = 11957.429
$ws.Range("M70").Value = -11687.429
$ws.Range("H73").Value = 12038.723
$ws.Range("I73").Value = 11957.429
$ws.Range("K73").Value = 11957.429
$ws.Range("M73").Value = -11021.429
$ws.Range("H81").Value = 46665.5
$ws.Range("J81").Value = 46665.5
$ws.Range("L81").Value = 46665.5
$ws.Range("N81").Value = -48661.5
$ws.Range("H84").Value = 46665.5
$ws.Range("J84").Value = 46665.5
$ws.Range("L84").Value = 139996.5
$ws.Range("N84").Value = -149980.5
$ws.Range("H97").Value = 1323.36
$ws.Range("I97").Value = 1152.2632
$ws.Range("J97").Value = 1865.1666
$ws.Range("K97").Value = 1152.2632
$ws.Range("L97").Value = 1865.1666
$ws.Range("M97").Value = -656.2632000000001
$ws.Range("N97").Value = -2857.1666
$ws.Range("H102").Value = 2754.8445
$ws.Range("I102").Value = 1938.5454
$ws.Range("K102").Value = 1938.5454
$ws.Range("M102").Value = -316.5454
$ws.Range("H122").Value = 5987.2925
$ws.Range("I122").Value = 6143.4644
$ws.Range("J122").Value = 5650.923
$ws.Range("K122").Value = 18430.3932
$ws.Range("L122").Value = 16952.769
$ws.Range("M122").Value = -15980.3932
$ws.Range("N122").Value = -21852.769
$ws.Range("H132").Value = 3356.818
$ws.Range("I132").Value = 1602.4
$ws.Range("J132").Value = 4818.8335
$ws.Range("K132").Value = 4807.200000000001
$ws.Range("L132").Value = 14456.5005
$ws.Range("M132").Value = -2277.200000000001
$ws.Range("N132").Value = -19516.5005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 14215.333
$ws.Range("I7").Value = 7058.4
$ws.Range("K7").Value = 7058.4
$ws.Range("M7").Value = -6946.4
$ws.Range("H16").Value = 2135.6365
$ws.Range("I16").Value = 1574.875
$ws.Range("K16").Value = 1574.875
$ws.Range("M16").Value = -1404.875
$ws.Range("H22").Value = 2821.5925
$ws.Range("I22").Value = 1593.0625
$ws.Range("J22").Value = 4608.5454
$ws.Range("K22").Value = 1593.0625
$ws.Range("L22").Value = 4608.5454
$ws.Range("M22").Value = -1298.0625
$ws.Range("N22").Value = -5198.5454
$ws.Range("H27").Value = 2821.5925
$ws.Range("I27").Value = 1593.0625
$ws.Range("J27").Value = 4608.5454
$ws.Range("K27").Value = 1593.0625
$ws.Range("L27").Value = 4608.5454
$ws.Range("M27").Value = -1486.0625
$ws.Range("N27").Value = -4822.5454
$ws.Range("H40").Value = 13810.5
$ws.Range("I40").Value = 15139.5
$ws.Range("J40").Value = 11152.5
$ws.Range("K40").Value = 15139.5
$ws.Range("L40").Value = 11152.5
$ws.Range("M40").Value = -15003.5
$ws.Range("N40").Value = -11424.5
$ws.Range("H46").Value = 3593.0833
$ws.Range("I46").Value = 2083.4167
$ws.Range("K46").Value = 2083.4167
$ws.Range("M46").Value = -1895.4167
$ws.Range("H50").Value = 44999
$ws.Range("I50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("M50").ClearContents()
$ws.Range("H53").Value = 5302
$ws.Range("I53").Value = 7003
$ws.Range("J53").Value = 1900
$ws.Range("K53").Value = 7003
$ws.Range("L53").Value = 1900
$ws.Range("M53").Value = -6485
$ws.Range("N53").Value = -2936
$ws.Range("H55").Value = 1713.6786
$ws.Range("I55").Value = 477.0625
$ws.Range("K55").Value = 477.0625
$ws.Range("M55").Value = -304.0625
$ws.Range("H57").Value = 39000
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 39000
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 39000
$ws.Range("M57").ClearContents()
$ws.Range("N57").Value = -40132
$ws.Range("H61").Value = 3578.721
$ws.Range("I61").Value = 2979.6
$ws.Range("J61").Value = 6199.875
$ws.Range("K61").Value = 2979.6
$ws.Range("L61").Value = 6199.875
$ws.Range("M61").Value = -2777.6
$ws.Range("N61").Value = -6603.875
$ws.Range("H68").Value = 4865.8335
$ws.Range("I68").Value = 4849.5
$ws.Range("J68").Value = 4874
$ws.Range("K68").Value = 4849.5
$ws.Range("L68").Value = 4874
$ws.Range("M68").Value = -4100.5
$ws.Range("N68").Value = -6372
$ws.Range("H71").Value = 4865.8335
$ws.Range("I71").Value = 4849.5
$ws.Range("J71").Value = 4874
$ws.Range("K71").Value = 24247.5
$ws.Range("L71").Value = 24370
$ws.Range("M71").Value = -20503.5
$ws.Range("N71").Value = -31858
$ws.Range("H87").Value = 50000
$ws.Range("I87").Value = 50000
$ws.Range("K87").Value = 50000
$ws.Range("M87").Value = -48877
$ws.Range("H90").Value = 50000
$ws.Range("I90").Value = 50000
$ws.Range("K90").Value = 150000
$ws.Range("M90").Value = -144384
$ws.Range("H93").Value = 1913.9615
$ws.Range("I93").Value = 1987
$ws.Range("J93").Value = 1715.7142
$ws.Range("K93").Value = 1987
$ws.Range("L93").Value = 1715.7142
$ws.Range("M93").Value = -739
$ws.Range("N93").Value = -4211.7142
$ws.Range("H96").Value = 49000
$ws.Range("J96").Value = 49000
$ws.Range("L96").Value = 49000
$ws.Range("N96").Value = -54492
$ws.Range("H113").Value = 3578.721
$ws.Range("I113").Value = 2979.6
$ws.Range("J113").Value = 6199.875
$ws.Range("K113").Value = 2979.6
$ws.Range("L113").Value = 6199.875
$ws.Range("M113").Value = -809.5999999999999
$ws.Range("N113").Value = -10539.875
$ws.Range("H126").Value = 14215.333
$ws.Range("I126").Value = 7058.4
$ws.Range("K126").Value = 21175.2
$ws.Range("M126").Value = -18705.2
$ws.Range("H127").Value = 88200
$ws.Range("J127").Value = 88200
$ws.Range("L127").Value = 88200
$ws.Range("N127").Value = -98120
$ws.Range("H132").Value = 6422.6665
$ws.Range("I132").Value = 5096.6665
$ws.Range("J132").Value = 7306.6665
$ws.Range("K132").Value = 15289.9995
$ws.Range("L132").Value = 21919.9995
$ws.Range("M132").Value = -12759.9995
$ws.Range("N132").Value = -26979.9995
$ws.Range("H136").Value = 8753.729499999999
$ws.Range("I136").Value = 6341.5
$ws.Range("J136").Value = 14455.363
$ws.Range("K136").Value = 19024.5
$ws.Range("L136").Value = 43366.089
$ws.Range("M136").Value = -16474.5
$ws.Range("N136").Value = -48466.089

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H53").Value = 1800
$ws.Range("J53").Value = 1800
$ws.Range("L53").Value = 1800
$ws.Range("N53").Value = -3014
$ws.Range("H62").Value = 12000.5
$ws.Range("J62").Value = 13800.6
$ws.Range("L62").Value = 13800.6
$ws.Range("N62").Value = -15048.6
$ws.Range("H65").Value = 12000.5
$ws.Range("J65").Value = 13800.6
$ws.Range("L65").Value = 69003
$ws.Range("N65").Value = -75243
$ws.Range("H81").Value = 4266.2856
$ws.Range("J81").Value = 4899.6
$ws.Range("L81").Value = 9799.200000000001
$ws.Range("N81").Value = -11921.2
$ws.Range("H82").Value = 49997.332
$ws.Range("J82").Value = 49997.332
$ws.Range("L82").Value = 49997.332
$ws.Range("N82").Value = -50763.332
$ws.Range("H84").Value = 4266.2856
$ws.Range("J84").Value = 4899.6
$ws.Range("L84").Value = 48996
$ws.Range("N84").Value = -59604
$ws.Range("H85").Value = 49997.332
$ws.Range("J85").Value = 49997.332
$ws.Range("L85").Value = 49997.332
$ws.Range("N85").Value = -52649.332
$ws.Range("H96").Value = 1443.1666
$ws.Range("I96").Value = 1131.8
$ws.Range("K96").Value = 1131.8
$ws.Range("M96").Value = 241.2
$ws.Range("H122").Value = 2606.1924
$ws.Range("I122").Value = 1814
$ws.Range("K122").Value = 5442
$ws.Range("M122").Value = -2992
$ws.Range("H126").Value = 3554.889
$ws.Range("I126").Value = 2399.6
$ws.Range("J126").Value = 4999
$ws.Range("K126").Value = 7198.799999999999
$ws.Range("L126").Value = 14997
$ws.Range("M126").Value = -4728.799999999999
$ws.Range("N126").Value = -19937
$ws.Range("H132").Value = 17669.666
$ws.Range("I132").Value = 17002
$ws.Range("K132").Value = 51006
$ws.Range("M132").Value = -48476
$ws.Range("H136").Value = 3883.4736
$ws.Range("I136").Value = 1985.5333
$ws.Range("J136").Value = 11000.75
$ws.Range("K136").Value = 5956.5999
$ws.Range("L136").Value = 33002.25
$ws.Range("M136").Value = -3406.5999
$ws.Range("N136").Value = -38102.25

Write-Host "Applied all changes"